# Auto-generated script applying the market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 246.75
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 1997.8148
$ws.Range("I86").Value = 1121.8
$ws.Range("J86").Value = 3092.8333
$ws.Range("K86").Value = 1121.8
$ws.Range("L86").Value = 3092.8333
$ws.Range("M86").Value = 1.200000000000045
$ws.Range("N86").Value = -5338.8333
$ws.Range("H89").Value = 1997.8148
$ws.Range("I89").Value = 1121.8
$ws.Range("J89").Value = 3092.8333
$ws.Range("K89").Value = 5609
$ws.Range("L89").Value = 15464.1665
$ws.Range("M89").Value = 7
$ws.Range("N89").Value = -26696.1665
$ws.Range("H98").Value = 2037.5625
$ws.Range("I98").Value = 1673
$ws.Range("K98").Value = 1673
$ws.Range("M98").Value = -175
$ws.Range("H103").Value = 506.1905
$ws.Range("J103").Value = 699
$ws.Range("L103").Value = 2097
$ws.Range("N103").Value = -3269
$ws.Range("H107").Value = 556.6896400000001
$ws.Range("I107").Value = 573.61536
$ws.Range("K107").Value = 573.61536
$ws.Range("M107").Value = 1346.38464
$ws.Range("H122").Value = 2037.5625
$ws.Range("I122").Value = 1673
$ws.Range("K122").Value = 5019
$ws.Range("M122").Value = -2569
$ws.Range("H127").Value = 745.5
$ws.Range("I127").Value = 745.5
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2236.5
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 2723.5
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1953.081
$ws.Range("I132").Value = 1522.9697
$ws.Range("K132").Value = 4568.909100000001
$ws.Range("M132").Value = -2038.909100000001
$ws.Range("H137").Value = 55557732
$ws.Range("I137").Value = 100001060
$ws.Range("K137").Value = 300003180
$ws.Range("M137").Value = -300000630
$ws.Range("H138").Value = 4874.075
$ws.Range("I138").Value = 2979.8
$ws.Range("J138").Value = 5144.6855
$ws.Range("K138").Value = 8939.400000000001
$ws.Range("L138").Value = 15434.0565
$ws.Range("M138").Value = -3799.400000000001
$ws.Range("N138").Value = -25714.0565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5728
$ws.Range("I32").Value = 4234.9165
$ws.Range("K32").Value = 4234.9165
$ws.Range("M32").Value = -3947.9165
$ws.Range("H88").Value = 1152.7273
$ws.Range("I88").Value = 1251
$ws.Range("J88").Value = 1054.4546
$ws.Range("K88").Value = 1251
$ws.Range("L88").Value = 1054.4546
$ws.Range("M88").Value = -845
$ws.Range("N88").Value = -1866.4546
$ws.Range("H91").Value = 1152.7273
$ws.Range("I91").Value = 1251
$ws.Range("J91").Value = 1054.4546
$ws.Range("K91").Value = 1251
$ws.Range("L91").Value = 1054.4546
$ws.Range("M91").Value = 153
$ws.Range("N91").Value = -3862.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 58711.5
$ws.Range("J95").Value = 58711.5
$ws.Range("L95").Value = 58711.5
$ws.Range("N95").Value = -64203.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6732.364
$ws.Range("I31").Value = 6194.048
$ws.Range("J31").Value = 7674.4165
$ws.Range("K31").Value = 6194.048
$ws.Range("L31").Value = 7674.4165
$ws.Range("M31").Value = -5899.048
$ws.Range("N31").Value = -8264.416499999999
$ws.Range("H34").Value = 6732.364
$ws.Range("I34").Value = 6194.048
$ws.Range("J34").Value = 7674.4165
$ws.Range("K34").Value = 6194.048
$ws.Range("L34").Value = 7674.4165
$ws.Range("M34").Value = -5992.048
$ws.Range("N34").Value = -8078.4165
$ws.Range("H35").Value = 10721.111
$ws.Range("I35").Value = 10721.111
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 10721.111
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -10427.111
$ws.Range("N35").ClearContents()
$ws.Range("H105").Value = 23277
$ws.Range("I105").Value = 32849.09
$ws.Range("J105").Value = 2218.4
$ws.Range("K105").Value = 32849.09
$ws.Range("L105").Value = 2218.4
$ws.Range("M105").Value = -31102.09
$ws.Range("N105").Value = -5712.4
$ws.Range("H107").Value = 627.82355
$ws.Range("I107").Value = 535.8182
$ws.Range("K107").Value = 535.8182
$ws.Range("M107").Value = 1384.1818
$ws.Range("H122").Value = 2392.7
$ws.Range("J122").Value = 2578.375
$ws.Range("L122").Value = 7735.125
$ws.Range("N122").Value = -12635.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 389.375
$ws.Range("J23").Value = 323.2
$ws.Range("L23").Value = 969.5999999999999
$ws.Range("N23").Value = -1439.6
$ws.Range("H34").Value = 740.7857
$ws.Range("I34").Value = 811.8333
$ws.Range("J34").Value = 687.5
$ws.Range("K34").Value = 2435.4999
$ws.Range("L34").Value = 2062.5
$ws.Range("M34").Value = -2351.4999
$ws.Range("N34").Value = -2230.5
$ws.Range("H39").Value = 4166.8335
$ws.Range("I39").Value = 3998
$ws.Range("J39").Value = 4200.6
$ws.Range("K39").Value = 11994
$ws.Range("L39").Value = 12601.8
$ws.Range("M39").Value = -11700
$ws.Range("N39").Value = -13189.8
$ws.Range("H49").Value = 500
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H55").Value = 912.1429000000001
$ws.Range("I55").Value = 995
$ws.Range("J55").Value = 898.3333
$ws.Range("K55").Value = 2985
$ws.Range("L55").Value = 2694.9999
$ws.Range("M55").Value = -2808
$ws.Range("N55").Value = -3048.9999
$ws.Range("H88").Value = 4004.6667
$ws.Range("J88").Value = 4500
$ws.Range("L88").Value = 13500
$ws.Range("N88").Value = -14356
$ws.Range("H91").Value = 4004.6667
$ws.Range("J91").Value = 4500
$ws.Range("L91").Value = 13500
$ws.Range("N91").Value = -16464
$ws.Range("H131").Value = 1554.5897
$ws.Range("J131").Value = 1670.5758
$ws.Range("L131").Value = 5011.7274
$ws.Range("N131").Value = -15091.7274
$ws.Range("H140").Value = 2775.8438
$ws.Range("I140").Value = 2037.4783
$ws.Range("J140").Value = 4662.778
$ws.Range("K140").Value = 6112.4349
$ws.Range("L140").Value = 13988.334
$ws.Range("M140").Value = -932.4349000000002
$ws.Range("N140").Value = -24348.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 56491.95
$ws.Range("I122").Value = 78012.46000000001
$ws.Range("K122").Value = 234037.38
$ws.Range("M122").Value = -231587.38

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4383.4814
$ws.Range("I16").Value = 5612.65
$ws.Range("J16").Value = 871.5714
$ws.Range("K16").Value = 5612.65
$ws.Range("L16").Value = 871.5714
$ws.Range("M16").Value = -5442.65
$ws.Range("N16").Value = -1211.5714
$ws.Range("H82").Value = 723
$ws.Range("I82").Value = 651
$ws.Range("J82").Value = 975
$ws.Range("K82").Value = 651
$ws.Range("L82").Value = 975
$ws.Range("M82").Value = -290
$ws.Range("N82").Value = -1697
$ws.Range("H85").Value = 723
$ws.Range("I85").Value = 651
$ws.Range("J85").Value = 975
$ws.Range("K85").Value = 651
$ws.Range("L85").Value = 975
$ws.Range("M85").Value = 597
$ws.Range("N85").Value = -3471

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 260342.72
$ws.Range("I62").Value = 362559.8
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 362559.8
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -361935.8
$ws.Range("N62").Value = -6048
$ws.Range("H65").Value = 260342.72
$ws.Range("I65").Value = 362559.8
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 1812799
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -1809679
$ws.Range("N65").Value = -30240
$ws.Range("H87").Value = 99999
$ws.Range("J87").Value = 99999
$ws.Range("L87").Value = 99999
$ws.Range("N87").Value = -102495
$ws.Range("H90").Value = 99999
$ws.Range("J90").Value = 99999
$ws.Range("L90").Value = 299997
$ws.Range("N90").Value = -312477
$ws.Range("H100").Value = 1294.6666
$ws.Range("I100").Value = 1700
$ws.Range("J100").Value = 1092
$ws.Range("K100").Value = 3400
$ws.Range("L100").Value = 2184
$ws.Range("M100").Value = -2859
$ws.Range("N100").Value = -3266
$ws.Range("H126").Value = 1387.4706
$ws.Range("I126").Value = 1453.1072
$ws.Range("K126").Value = 4359.321599999999
$ws.Range("M126").Value = -1889.321599999999
